$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.490.24'
$ws.Range("E2").Value = '  -2.60%  '
$ws.Range("D3").Value = '1.805.17'
$ws.Range("E3").Value = '  -2.51%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.008'
$ws.Range("E4").Value = '  +0.70%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.007'
$ws.Range("E5").Value = '  +0.60%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '308.05'
$ws.Range("E6").Value = '  -1.66%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4559'
$ws.Range("E7").Value = '  -1.28%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3661'
$ws.Range("E8").Value = '  -1.36%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07125'
$ws.Range("E9").Value = '  -2.00%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8779'
$ws.Range("E10").Value = '  -0.70%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07738'
$ws.Range("E11").Value = '  -0.81%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '19.34'
$ws.Range("E12").Value = '  -3.12%  '
$ws.Range("D13").Value = '1.804.23'
$ws.Range("E13").Value = '  -6.93%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.271'
$ws.Range("E14").Value = '  -1.74%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.359'
$ws.Range("E15").Value = '  -2.21%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '86.01'
$ws.Range("E16").Value = '  -5.59%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.008'
$ws.Range("E17").Value = '  +0.70%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008575'
$ws.Range("E18").Value = '  -3.56%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.007'
$ws.Range("E19").Value = '  +0.53%  '
$ws.Range("D20").Value = '26.547.74'
$ws.Range("E20").Value = '  -2.43%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.24'
$ws.Range("E21").Value = '  -2.93%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.979'
$ws.Range("E22").Value = '  -1.50%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.43'
$ws.Range("E23").Value = '  -0.44%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.982'
$ws.Range("E24").Value = '  +1.47%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '150.87'
$ws.Range("E25").Value = '  -0.73%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '17.93'
$ws.Range("E26").Value = '  -2.74%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.030'
$ws.Range("E27").Value = '  -0.86%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '112.57'
$ws.Range("E28").Value = '  -2.52%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.834'
$ws.Range("E29").Value = '  -4.35%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.08661'
$ws.Range("E30").Value = '  -1.63%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.039'
$ws.Range("E31").Value = '  -2.54%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7291'
$ws.Range("E32").Value = '  -4.52%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.439'
$ws.Range("E33").Value = '  -1.22%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.113'
$ws.Range("E34").Value = '  -4.68%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.007'
$ws.Range("E35").Value = '  +0.64%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.536'
$ws.Range("E36").Value = '  -6.99%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.078'
$ws.Range("E37").Value = '  -0.10%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01931'
$ws.Range("E38").Value = '  -0.35%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05104'
$ws.Range("E39").Value = '  -2.41%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.885'
$ws.Range("E40").Value = '  -1.82%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.945'
$ws.Range("E41").Value = '  -1.50%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4995'
$ws.Range("E42").Value = '  -1.97%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1566'
$ws.Range("E43").Value = '  -3.59%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.123'
$ws.Range("E44").Value = '  -3.00%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.007'
$ws.Range("E45").Value = '  +0.65%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4600'
$ws.Range("E46").Value = '  -3.70%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '101.64'
$ws.Range("E47").Value = '  -1.04%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.914'
$ws.Range("E48").Value = '  -3.30%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.591'
$ws.Range("E49").Value = '  -2.55%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05986'
$ws.Range("E50").Value = '  -3.60%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '63.93'
$ws.Range("E51").Value = '  -2.72%  '
